# Update legacy GSC export data ("Chart" sheet):
#  - Drop the 5 oldest days (2025-08-26 .. 2025-08-30) by deleting rows 2-6,
#    which shifts every remaining row up by 5.
#  - Append 4 new days (2025-11-25 .. 2025-11-28) at the bottom with their
#    HTTPS/Non-HTTPS counts (0 for both, matching the existing trailing data).
#  - The net effect is the 91-row window becomes a 90-row window
#    (dimension A1:C92 -> A1:C91).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Drop the oldest 5 rows (2025-08-26 .. 2025-08-30); this shifts all
# subsequent rows up by 5 rows automatically.
$ws.Range("A2:A6").EntireRow.Delete()

# After the delete, the last populated row (previously row 92, date
# 2025-11-24) is now row 87. Append the 4 new trailing days after it.
$newDates = @("2025-11-25", "2025-11-26", "2025-11-27", "2025-11-28")
$startRow = 88

# Pre-format the new date cells as Text so the "yyyy-mm-dd"-looking
# strings are stored verbatim (matching every other Date cell in the
# column) instead of being auto-converted into a date serial value.
$ws.Range("A88:A91").NumberFormat = "@"

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newDates[$i]
    $ws.Cells.Item($r, 2).Value = 0.0
    $ws.Cells.Item($r, 3).Value = 0.0
}
